# "demeaned rs scores (binning)"
# Composite_Reward and Composite_Reward_Squared columns are replaced with
# demeaned/binned placeholder values -- every row's B (Composite_Reward) and
# C (Composite_Reward_Squared) value becomes 0. The Subject id column (A)
# and the header row (A1:C1) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:C51").Value = 0
